# Scheduled-runner refresh of cached market-price / profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the Malboro_Profits
# sheets. Values are plain numeric snapshots (no formulas involved), so
# each changed cell is just re-written with its new cached figure. A
# couple of rows (LTW!L68:M68, LTW!L71:M71) had their HQ-price/profit
# cells removed entirely in the source edit (no longer quoted), so those
# are cleared instead of zero-filled.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1200
$ws.Cells.Item(40, 9).Value = 700
$ws.Cells.Item(40, 11).Value = 700
$ws.Cells.Item(40, 13).Value = -525
$ws.Cells.Item(46, 8).Value = 2095
$ws.Cells.Item(46, 9).Value = 1734
$ws.Cells.Item(46, 10).Value = 2456
$ws.Cells.Item(46, 11).Value = 5202
$ws.Cells.Item(46, 12).Value = 7368
$ws.Cells.Item(46, 13).Value = -5083
$ws.Cells.Item(46, 14).Value = -7606
$ws.Cells.Item(60, 8).Value = 2095
$ws.Cells.Item(60, 9).Value = 1734
$ws.Cells.Item(60, 10).Value = 2456
$ws.Cells.Item(60, 11).Value = 5202
$ws.Cells.Item(60, 12).Value = 7368
$ws.Cells.Item(60, 13).Value = -4718
$ws.Cells.Item(60, 14).Value = -8336
$ws.Cells.Item(125, 8).Value = 4136471.2
$ws.Cells.Item(125, 9).Value = 5055290
$ws.Cells.Item(125, 10).Value = 1787
$ws.Cells.Item(125, 11).Value = 45497610
$ws.Cells.Item(125, 12).Value = 16083
$ws.Cells.Item(125, 13).Value = -45495150
$ws.Cells.Item(125, 14).Value = -21003
$ws.Cells.Item(137, 8).Value = 6000.9775
$ws.Cells.Item(137, 9).Value = 1395.8108
$ws.Cells.Item(137, 11).Value = 4187.4324
$ws.Cells.Item(137, 13).Value = -1637.4324

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4632258
$ws.Cells.Item(32, 9).Value = 4903822
$ws.Cells.Item(32, 10).Value = 15665.667
$ws.Cells.Item(32, 11).Value = 4903822
$ws.Cells.Item(32, 12).Value = 15665.667
$ws.Cells.Item(32, 13).Value = -4903535
$ws.Cells.Item(32, 14).Value = -16239.667
$ws.Cells.Item(74, 8).Value = 21567.139
$ws.Cells.Item(74, 9).Value = 1474.2609
$ws.Cells.Item(74, 10).Value = 98589.836
$ws.Cells.Item(74, 11).Value = 1474.2609
$ws.Cells.Item(74, 12).Value = 98589.836
$ws.Cells.Item(74, 13).Value = -600.2609
$ws.Cells.Item(74, 14).Value = -100337.836
$ws.Cells.Item(77, 8).Value = 21567.139
$ws.Cells.Item(77, 9).Value = 1474.2609
$ws.Cells.Item(77, 10).Value = 98589.836
$ws.Cells.Item(77, 11).Value = 7371.3045
$ws.Cells.Item(77, 12).Value = 492949.18
$ws.Cells.Item(77, 13).Value = -3003.3045
$ws.Cells.Item(77, 14).Value = -501685.18
$ws.Cells.Item(132, 8).Value = 3685317
$ws.Cells.Item(132, 9).Value = 3123.76
$ws.Cells.Item(132, 11).Value = 9371.280000000001
$ws.Cells.Item(132, 13).Value = -6841.280000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 17369.79
$ws.Cells.Item(31, 9).Value = 1362.3636
$ws.Cells.Item(31, 11).Value = 1362.3636
$ws.Cells.Item(31, 13).Value = -1067.3636
$ws.Cells.Item(34, 8).Value = 17369.79
$ws.Cells.Item(34, 9).Value = 1362.3636
$ws.Cells.Item(34, 11).Value = 1362.3636
$ws.Cells.Item(34, 13).Value = -1160.3636
$ws.Cells.Item(99, 8).Value = 4728595
$ws.Cells.Item(99, 9).Value = 2970010.2
$ws.Cells.Item(99, 10).Value = 10004350
$ws.Cells.Item(99, 11).Value = 2970010.2
$ws.Cells.Item(99, 12).Value = 10004350
$ws.Cells.Item(99, 13).Value = -2968512.2
$ws.Cells.Item(99, 14).Value = -10007346
$ws.Cells.Item(126, 8).Value = 4728595
$ws.Cells.Item(126, 9).Value = 2970010.2
$ws.Cells.Item(126, 10).Value = 10004350
$ws.Cells.Item(126, 11).Value = 8910030.600000001
$ws.Cells.Item(126, 12).Value = 30013050
$ws.Cells.Item(126, 13).Value = -8907560.600000001
$ws.Cells.Item(126, 14).Value = -30017990
$ws.Cells.Item(132, 8).Value = 114383624
$ws.Cells.Item(132, 9).Value = 4955.25
$ws.Cells.Item(132, 11).Value = 14865.75
$ws.Cells.Item(132, 13).Value = -12335.75
$ws.Cells.Item(134, 8).Value = 41674630
$ws.Cells.Item(134, 9).Value = 2070.818
$ws.Cells.Item(134, 10).Value = 76936030
$ws.Cells.Item(134, 11).Value = 6212.454000000001
$ws.Cells.Item(134, 12).Value = 230808090
$ws.Cells.Item(134, 13).Value = -3677.454000000001
$ws.Cells.Item(134, 14).Value = -230813160
$ws.Cells.Item(141, 8).Value = 366000
$ws.Cells.Item(141, 9).Value = 39000
$ws.Cells.Item(141, 10).Value = 529500
$ws.Cells.Item(141, 11).Value = 39000
$ws.Cells.Item(141, 12).Value = 529500
$ws.Cells.Item(141, 13).Value = -33820
$ws.Cells.Item(141, 14).Value = -539860

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 85.58064
$ws.Cells.Item(2, 9).Value = 107
$ws.Cells.Item(2, 10).Value = 46.636364
$ws.Cells.Item(2, 11).Value = 642
$ws.Cells.Item(2, 12).Value = 279.818184
$ws.Cells.Item(2, 13).Value = -529
$ws.Cells.Item(2, 14).Value = -505.818184
$ws.Cells.Item(7, 8).Value = 153.73334
$ws.Cells.Item(7, 9).Value = 124
$ws.Cells.Item(7, 11).Value = 372
$ws.Cells.Item(7, 13).Value = -260
$ws.Cells.Item(23, 8).Value = 372.33334
$ws.Cells.Item(23, 9).Value = 60
$ws.Cells.Item(23, 10).Value = 434.8
$ws.Cells.Item(23, 11).Value = 180
$ws.Cells.Item(23, 12).Value = 1304.4
$ws.Cells.Item(23, 13).Value = 55
$ws.Cells.Item(23, 14).Value = -1774.4
$ws.Cells.Item(113, 8).Value = 1150.7084
$ws.Cells.Item(113, 10).Value = 1166.75
$ws.Cells.Item(113, 12).Value = 3500.25
$ws.Cells.Item(113, 14).Value = -7840.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1098856.6
$ws.Cells.Item(122, 9).Value = 1419067
$ws.Cells.Item(122, 10).Value = 992.5714
$ws.Cells.Item(122, 11).Value = 4257201
$ws.Cells.Item(122, 12).Value = 2977.7142
$ws.Cells.Item(122, 13).Value = -4254751
$ws.Cells.Item(122, 14).Value = -7877.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(8, 8).Value = 179168
$ws.Cells.Item(8, 10).Value = 179168
$ws.Cells.Item(8, 12).Value = 179168
$ws.Cells.Item(8, 14).Value = -179448
$ws.Cells.Item(16, 8).Value = 76925944
$ws.Cells.Item(16, 9).Value = 100002670
$ws.Cells.Item(16, 11).Value = 100002670
$ws.Cells.Item(16, 13).Value = -100002500
$ws.Cells.Item(40, 8).Value = 1594290.5
$ws.Cells.Item(40, 10).Value = 3466003.2
$ws.Cells.Item(40, 12).Value = 3466003.2
$ws.Cells.Item(40, 14).Value = -3466275.2
$ws.Cells.Item(55, 8).Value = 1456.6786
$ws.Cells.Item(55, 9).Value = 982
$ws.Cells.Item(55, 10).Value = 2458.7778
$ws.Cells.Item(55, 11).Value = 982
$ws.Cells.Item(55, 12).Value = 2458.7778
$ws.Cells.Item(55, 13).Value = -809
$ws.Cells.Item(55, 14).Value = -2804.7778
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).ClearContents()
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(68, 14).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).ClearContents()
$ws.Cells.Item(71, 13).ClearContents()
$ws.Cells.Item(71, 14).Value = 0
$ws.Cells.Item(132, 8).Value = 1629967.1
$ws.Cells.Item(132, 9).Value = 4244.926
$ws.Cells.Item(132, 10).Value = 4373373.5
$ws.Cells.Item(132, 11).Value = 12734.778
$ws.Cells.Item(132, 12).Value = 13120120.5
$ws.Cells.Item(132, 13).Value = -10204.778
$ws.Cells.Item(132, 14).Value = -13125180.5
$ws.Cells.Item(136, 8).Value = 1586250.9
$ws.Cells.Item(136, 9).Value = 52299.75
$ws.Cells.Item(136, 10).Value = 1995304.5
$ws.Cells.Item(136, 11).Value = 156899.25
$ws.Cells.Item(136, 12).Value = 5985913.5
$ws.Cells.Item(136, 13).Value = -154349.25
$ws.Cells.Item(136, 14).Value = -5991013.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 148795.67
$ws.Cells.Item(62, 10).Value = 161145.12
$ws.Cells.Item(62, 12).Value = 161145.12
$ws.Cells.Item(62, 14).Value = -162393.12
$ws.Cells.Item(65, 8).Value = 148795.67
$ws.Cells.Item(65, 10).Value = 161145.12
$ws.Cells.Item(65, 12).Value = 805725.6
$ws.Cells.Item(65, 14).Value = -811965.6
$ws.Cells.Item(136, 8).Value = 321505.28
$ws.Cells.Item(136, 10).Value = 660656.6
$ws.Cells.Item(136, 12).Value = 1981969.8
$ws.Cells.Item(136, 14).Value = -1987069.8
